$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16:B21").Value = "Done"

$fc = $ws.Cells.FormatConditions
$fc.Delete()
